$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Tipo", matching the style/format of the existing header row
$ws.Range("D1").Value = "Tipo"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update MSE (column B) and R2 (column C) values for rows 2-4
$ws.Range("B2:B4").Value = 0.3968241615722563
$ws.Range("C2:C4").Value = 0.9944915471856128

# Add the new "Tipo" column values for rows 2-4
$ws.Range("D2:D4").Value = "single"
